# Auto update stock data
# Updates the "latest" date row for each ticker block from 2025/11/06 -> 2025/11/07
# and refreshes the accompanying metric values, matching the daily data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param(
        [string]$Address,
        [string]$Value
    )
    # Force the cell to be treated as text so Excel doesn't auto-convert
    # date-like or numeric-like strings into dates/numbers.
    $rng = $ws.Range($Address)
    $rng.NumberFormat = "@"
    $rng.Value = $Value
}

# Row 2 (AA)
Set-TextCell "A2" "2025/11/07"
Set-TextCell "B2" "4.47"

# Row 8 (RIO)
Set-TextCell "A8" "2025/11/07"
Set-TextCell "B8" "7.42"

# Row 14 (NHY)
Set-TextCell "A14" "2025/11/07"
Set-TextCell "B14" "2.71"

# Row 20 (RS)
Set-TextCell "A20" "2025/11/07"
Set-TextCell "B20" "12.18"

# Row 26 (KALU)
Set-TextCell "A26" "2025/11/07"
Set-TextCell "B26" "9.71"

# Row 32 (RYI)
Set-TextCell "A32" "2025/11/07"
Set-TextCell "B32" "24.80"

# Row 38 (BVB:ALR) - only date changes
Set-TextCell "A38" "2025/11/07"

# Row 44 (ULTR)
Set-TextCell "A44" "2025/11/07"
Set-TextCell "B44" "11.23"

# Row 50 (BHE) - date, EBITDA, Debt/Equity, and Altman Z-Score (now populated)
Set-TextCell "A50" "2025/11/07"
Set-TextCell "B50" "11.72"
Set-TextCell "C50" "0.30"
$ws.Range("G50").Value = 3.22

# Rows 51-55 (BHE historical rows) - Altman Z-Score now populated
$ws.Range("G51").Value = 3.22
$ws.Range("G52").Value = 3.22
$ws.Range("G53").Value = 3.22
$ws.Range("G54").Value = 3.22
$ws.Range("G55").Value = 3.22

# Row 56 (CLS) - only date changes
Set-TextCell "A56" "2025/11/07"

# Row 62 (JABIL)
Set-TextCell "A62" "2025/11/07"
Set-TextCell "B62" "11.67"

# Row 68 (FLEX)
Set-TextCell "A68" "2025/11/07"
Set-TextCell "B68" "13.17"

# Row 74 (MKS)
Set-TextCell "A74" "2025/11/07"
Set-TextCell "B74" "15.94"
